$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("RangeContains")
$ws2 = $wb.Worksheets.Item("Contains")

# Update float/Float/FloatValue sample values from 1.05 to 1.1 on both sheets.
# Prefix with a leading apostrophe so Excel stores these as plain text
# (shared strings) instead of interpreting the leading "=" as a formula.
# Order matters for shared-string table append order: float, then
# FloatValue, then Float.
$ws1.Range("C12").Value = "'= (float) 1.1"
$ws1.Range("C26").Value = "'= (FloatValue) 1.1"
$ws1.Range("C18").Value = "'= (Float) 1.1"

$ws2.Range("C13").Value = "'= (float) 1.1"
$ws2.Range("C27").Value = "'= (FloatValue) 1.1"
$ws2.Range("C19").Value = "'= (Float) 1.1"

# Update the selection on each sheet and make "RangeContains" the active
# (selected) tab again, matching the saved view state.
$ws2.Activate()
$ws2.Range("C20").Select()

$ws1.Activate()
$ws1.Range("C27").Select()
